$p = $ppt.ActivePresentation

for ($i = 1; $i -le 6; $i++) {
    $s = $p.Slides.Item($i)
    $null = $s.TimeLine.MainSequence
}
